$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: delete B2, D2, E2 and tweak C2
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -5.0379295087167808
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 tweaks
$ws.Range("B3").Value = -6.5016201590062561
$ws.Range("C3").Value = 8.8332508674856403
$ws.Range("D3").Value = -0.79204125872306064
$ws.Range("E3").Value = 27.321245754272983

# Update selection to match new state
$ws.Range("B1:E3").Select()
